$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 60")

# Resize the textbox (TextBox 60)
$shp.Width = 9520765 / 12700
$shp.Height = 5078313 / 12700

$tr = $shp.TextFrame.TextRange

# Merge the "Because of parse tree complexity..." paragraph (lvl=1 bullet)
# into the preceding "No unit testing..." paragraph, then rewrite the text.
$para2 = $tr.Paragraphs(2)
$para3 = $tr.Paragraphs(3)
$para3.Delete()

$para2.Text = "No unit testing for main code because Black Box tests cover the majority of use "
$para2.InsertAfter("cases, along with cases we wouldn’t have come up with") | Out-Null
